$d = $word.ActiveDocument

# 1. CarType row: dif_branch_return_price -> dif_branch_ret_price
$d.Content.Find.Execute("dif_branch_return_price", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "dif_branch_ret_price", 2)

# 2. Car row: drop the trailing ", branch_id" column
$d.Content.Find.Execute(", make, model, year, colour, license_plate, type, branch_id", $true, $false, $false, $false, $false, `
                         $true, 1, $false, ", make, model, year, colour, license_plate, type", 2)

# 3. Drop "branch_id references Branch." from the end of the "type references CarType" note,
#    keeping the trailing comma (and the line break that follows it untouched).
$d.Content.Find.Execute(", branch_id references Branch.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, ",", 2)

# 4. Rewrite the Rental row definition with the new column list / order. Scoped to the
#    Rental paragraph and to the text *after* "reservation_id, " so the underline
#    formatting on the reservation_id primary key is left untouched.
$pRental = $d.Paragraphs.Item(11)
$pRental.Range.Find.Execute("from_date, to_date, price(), customer_id, renting_branch_id, returning_branch_id, vin", `
                             $true, $false, $false, $false, $false, `
                             $true, 1, $false, `
                             "from_date, to_date, customer_id, vin, branch_id_pickup, branch_id_return, price()", 2)

# 5. Update the three Rental foreign-key notes (scoped per-paragraph so the repeated
#    "vin references Car" phrase can't bleed across lines).
$pRenting = $d.Paragraphs.Item(13)
$pRenting.Range.Find.Execute("renting_branch_id references Branch", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "vin references Car", 2)

$pReturning = $d.Paragraphs.Item(14)
$pReturning.Range.Find.Execute("returning_branch_id references Branch", $true, $false, $false, $false, $false, `
                                $true, 1, $false, "branch_id_pickup references Branch", 2)

$pVin = $d.Paragraphs.Item(15)
$pVin.Range.Find.Execute("vin references Car", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "branch_id_return references Branch", 2)

# 6. Append the two new login tables (with a blank paragraph between them) after the
#    last paragraph of the document.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$pLogin1 = $d.Paragraphs.Item($d.Paragraphs.Count)
$pLogin1.Style = "Normal"
$pLogin1.Range.Text = "CustomerLogin (customer_id, username, password)"

$pLogin1.Range.InsertParagraphAfter()
$pBlank = $d.Paragraphs.Item($d.Paragraphs.Count)
$pBlank.Style = "Normal"

$pBlank.Range.InsertParagraphAfter()
$pLogin2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$pLogin2.Style = "Normal"
$pLogin2.Range.Text = "EmployeeLogin (employee_id, username, password)"
